# Update of all values to match PDF edition 10 (commit 1)
$wb = $excel.ActiveWorkbook

$wsAdditional = $wb.Worksheets.Item("additional")

# Switch to the "additional" sheet before editing it
$wsAdditional.Activate()

# Add the new 2022 column of data (Year/Minutes table grows from H to I)
$wsAdditional.Range("I2").Value = 2022
$wsAdditional.Range("I3").Value = 3.22

# Leave the selection where the author left it when saving
$wsAdditional.Range("G8").Select()

$wb.Save()
